# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# per the latest scraped data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 190
$wsExpo.Range("F3").Value = 512
$wsExpo.Range("F4").Value = 33
$wsExpo.Range("F5").Value = 23
$wsExpo.Range("F6").Value = 14
$wsExpo.Range("F7").Value = 32
$wsExpo.Range("F8").Value = 24
$wsExpo.Range("F9").Value = 174
$wsExpo.Range("F10").Value = 2511
$wsExpo.Range("F11").Value = 16

# --- Sheet "全部类型" (fourth sheet) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 190
$wsAll.Range("F4").Value = 512
$wsAll.Range("F5").Value = 33
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 14
$wsAll.Range("F8").Value = 32
$wsAll.Range("F9").Value = 24
$wsAll.Range("F10").Value = 174
$wsAll.Range("F11").Value = 2511
$wsAll.Range("F12").Value = 16
